$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update today's blog cell (C2) with the new "ser" value (56 -> 61)
$ws.Range("C2").Value = "type: blog`nwidth: 2`nheight: 1`nser: 61"

# Update the active selection to reflect today's blog cell
[void]$ws.Range("C2").Select()
